# Update column G ("K") values in the save-data sheet.
# The commit message indicates the "K" column (strikes) was regenerated
# from a different/new source ("use K instead of Strike#") so the new
# numbers are written directly, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newK = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 2
    9  = 2
    10 = 2
    11 = 2
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    20 = 0
    21 = 1
    22 = 0
    23 = 2
    24 = 1
    25 = 1
    26 = 3
    27 = 1
    28 = 2
    29 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
